# Update the "Förändrad" date column (C) for rows 2-28.
# Each cell currently holds the serial date 45481 (2024-07-08) and must be
# bumped by one day to 45482 (2024-07-09), matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45481) {
        $cell.Value2 = 45482
    }
}
